$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bulk-updated
# from 45172 (2023-09-03) to 45175 (2023-09-06) for every data row (2-300).
$ws.Range("C2:C300").Value = 45175
